$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'41.527.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'2.468.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.97%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.34%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'313.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.51%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'91.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -2.24%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.14%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.33%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.515"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +3.16%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'32.53"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -3.08%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0791"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.81%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +0.28%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'2.850.03"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.99%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  -1.14%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'15.96"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +2.67%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'2.486.30"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.79%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -1.98%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'41.528.44"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.12%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +2.27%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.0₃0941"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.89%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'70.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.67%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'11.10"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.32%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'238.42"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.74%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.88%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +0.76%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -0.02%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'24.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.68%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'2.24"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.06%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -1.93%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'35.45"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -4.23%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'155.64"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.94%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'5.44"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.97%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -0.16%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.0759"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.47%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'17.25"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -3.87%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'2.37"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -4.22%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -5.92%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +2.96%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +0.70%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -4.73%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'4.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -2.84%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -0.52%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'1.945.51"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -2.31%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -0.92%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'18.72"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -5.69%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -3.45%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'9.07"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +2.73%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'2.709.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.06%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'97.24"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.06%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'67.24"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -3.09%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'52.33"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +2.90%  "
$ws.Range("E51").Style = "Normal"
